# Cronograma update: revised completion percentages (Terminacion %) for
# several activities, and move the active selection/viewport down a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G14").Value = 0.8
$ws.Range("G22").Value = 0.9
$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 0.8
$ws.Range("G25").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("G27").Value = 1

# Move selection to I33 (was I29) and scroll the view down so row 25 is
# at the top (was row 24).
$ws.Range("I33").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
